$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold a "Self?" flag indicating
# whether the row's country is the "self" country (Atlantis).
$ws.Columns("B:B").Insert()

# Header label for the newly inserted column
$ws.Range("B1").Value = "Self?"

# Column A still holds the country names (Atlantis, Brobdingnag, Carpania,
# Dinotopia, Erewhon in rows 2-6). Atlantis is the "self" country.
$ws.Range("B2").Value = "Yes"   # Atlantis
$ws.Range("B3").Value = "No"    # Brobdingnag
$ws.Range("B4").Value = "No"    # Carpania
$ws.Range("B5").Value = "No"    # Dinotopia
$ws.Range("B6").Value = "No"    # Erewhon

$ws.Range("C10").Select()
